# pushing recent version of EM - BH unbiased for pollock.
$wb = $excel.ActiveWorkbook

# --- Controls sheet: n_sims 100 -> 300 ---
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Range("B2").Value = 300

# --- Recruitment_Mortality sheet: M 1 -> 0.8 ---
$wsRecruit = $wb.Worksheets.Item("Recruitment_Mortality")
$wsRecruit.Range("B4").Value = 0.8

# Update selections to match the authored state
$wsRecruit.Activate()
$wsRecruit.Range("B5").Select() | Out-Null

# Controls becomes the active/selected sheet (tab selected), replacing
# Recruitment_Mortality which was active before.
$wsControls.Activate()
$wsControls.Range("B3").Select() | Out-Null
